# Update "想去人数" (want-to-go count) values in F column
# Sheet "展览" (index 1) and Sheet "全部类型" (index 4) both list the same
# events; bump the counts to the newly scraped totals.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item(1)   # 展览
$wsAllTypes   = $wb.Worksheets.Item(4)   # 全部类型

# Sheet 1 "展览" updates
$wsExhibition.Range("F3").Value  = 588
$wsExhibition.Range("F5").Value  = 103
$wsExhibition.Range("F7").Value  = 183
$wsExhibition.Range("F9").Value  = 1140
$wsExhibition.Range("F10").Value = 16487
$wsExhibition.Range("F15").Value = 197
$wsExhibition.Range("F17").Value = 11496
$wsExhibition.Range("F19").Value = 1144
$wsExhibition.Range("F20").Value = 4541
$wsExhibition.Range("F22").Value = 394

# Sheet 4 "全部类型" updates
$wsAllTypes.Range("F3").Value  = 588
$wsAllTypes.Range("F5").Value  = 103
$wsAllTypes.Range("F7").Value  = 183
$wsAllTypes.Range("F10").Value = 1140
$wsAllTypes.Range("F11").Value = 16487
$wsAllTypes.Range("F16").Value = 197
$wsAllTypes.Range("F20").Value = 11496
$wsAllTypes.Range("F22").Value = 1144
$wsAllTypes.Range("F23").Value = 4541
$wsAllTypes.Range("F25").Value = 394
